# Add the "11-07-2019 to 18-07-2019" summary sheet (same shape as the three
# existing daily-range sheets) and populate it with the per-client rollup
# built by iterating over the clients file.

$wb = $excel.ActiveWorkbook

# The workbook's original active sheet - used as the template for header
# formatting and as a safe source of already-"text" numeric-looking values
# (so new cells end up with the same shared-string/text storage the other
# three sheets use, instead of Excel's normal auto-number-detection).
$template = $wb.Worksheets.Item(1)

# New sheet goes after the last existing tab so it lands at the end (sheetId 4 / rId4).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "11-07-2019 to 18-07-2019"

# Client rollup rows: Clients / Enquires / Prices, alphabetically sorted,
# header first.
$clientNames = @(
  "All Drains LTD",
  "Asebestos",
  "Boiler Repair",
  "Conservatories",
  "Darrens Drains",
  "David",
  "Drain Division",
  "Drainage - UNSOLD",
  "Emergency Plumber - Unsold",
  "Essex SEO",
  "Evenings and Weekends",
  "Handyman - UNSOLD",
  "Manley Drainage LTD",
  "Mark Drainage",
  "North West Drain Service",
  "Ollie",
  "Pest Control",
  "Rapid Response Drain Care LTD",
  "Roofing",
  "SL Drainage",
  "Skip Hire",
  "Skip Hire 2",
  "Staffordshire Drainage Services ",
  "Staffordshire Drainage Services - Birmingham",
  "TC Drains",
  "UK Lie Detector Test",
  "USE Drainage Solutions",
  "Enquires"
)

# Matching "Prices" column per client, in the same order as $clientNames.
# Most clients settled at 450.0, a few (unsold-style rows) at 0.0, and two
# newly-seen clients priced at 112.5.
$priceKind = @(
  "450", "450", "450", "450", "450", "0", "450", "450", "450", "450",
  "450", "450", "450", "450", "450", "0", "450", "450", "450", "112",
  "112", "450", "450", "450", "450", "0", "450", "0"
)

# --- Header row: reuse sheet1's exact header formatting (bold/centered/bordered). ---
$template.Range("A1:C1").Copy()
$ws.Range("A1:C1").PasteSpecial(-4122)
$ws.Range("A1").Value = "Clients"
$ws.Range("B1").Value = "Enquires"
$ws.Range("C1").Value = "Prices"

for ($i = 0; $i -lt $clientNames.Length; $i++) {
  $row = $i + 2

  # Column A: client name - plain text, no special handling needed.
  $ws.Cells.Item($row, 1).Value = $clientNames[$i]

  # Column B: every row is "22.5" - copy the already-text "22.5" cell from
  # the template sheet so it stays text instead of becoming a number.
  $template.Range("B2").Copy()
  $ws.Cells.Item($row, 2).PasteSpecial(-4163)

  # Column C: "450.0" or "0.0" - both already exist as text elsewhere in the
  # workbook, so copy straight from there. Rows priced at "112.5" (a value
  # not seen anywhere else) are filled in afterwards.
  $kind = $priceKind[$i]
  if ($kind -eq "450") {
    $template.Range("C2").Copy()
    $ws.Cells.Item($row, 3).PasteSpecial(-4163)
  } elseif ($kind -eq "0") {
    $template.Range("C3").Copy()
    $ws.Cells.Item($row, 3).PasteSpecial(-4163)
  }
}

# --- A scratch cell used once to mint the literal text "112.5". A plain
# Value/Formula assignment of a numeric-looking string gets auto-detected
# as a number by Excel, so instead build it as a text-producing formula,
# then copy only the *computed value* (not the formula, not any format)
# into the real destination cells - that lands as a plain shared-string
# cell with no special number format, same as every other cell here. ---
$scratch = $ws.Range("Z1")
$scratch.Formula = '="112" & ".5"'

for ($i = 0; $i -lt $clientNames.Length; $i++) {
  if ($priceKind[$i] -eq "112") {
    $row = $i + 2
    $scratch.Copy()
    $ws.Cells.Item($row, 3).PasteSpecial(-4163)
  }
}

$scratch.Clear()

# Restore the original active sheet/selection so the new tab doesn't end up
# marked as the selected one.
$template.Activate()
$template.Range("A1").Select()
